# producao-de-jogos-digitais.xlsx - apply commit:
# "slides para definição dos recursos do jogo, definição da etapas e produtos"
#
# 1) analise-swot: rename competitor text + rewrite a threat/mitigation pair
# 2) add a brand-new "lista-mestra-de-recursos" worksheet with the feature
#    scoring table, legend, and sort order
# 3) tidy up view/selection state to match what Excel leaves behind after
#    the edit session

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. analise-swot (sheet1): text updates
# ---------------------------------------------------------------------------
$wsSwot = $wb.Worksheets.Item("analise-swot")
$wsSwot.Range("A2").Value = "O principal concorrente do Jogo da WHIP é o Jogo da Velha"
$wsSwot.Range("A7").Value = "o fator sorte pode ser determinante e prevalecer às habilidades do jogador."
$wsSwot.Range("B7").Value = "desenvolver um mecanismo que favoreça a melhor estratégia nas disputas por casas."

# ---------------------------------------------------------------------------
# 2. New worksheet: lista-mestra-de-recursos (added after the last sheet)
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsRec = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsRec.Name = "lista-mestra-de-recursos"

# column widths
$wsRec.Columns.Item(1).ColumnWidth = 13.28515625
$wsRec.Columns.Item(2).ColumnWidth = 67.42578125
$wsRec.Columns.Item(6).ColumnWidth = 10.85546875
$wsRec.Columns.Item(7).ColumnWidth = 10.85546875

# title banner (row 2), merged A2:H2, left aligned w/ bottom border
$wsRec.Range("A2:H2").Merge()
$wsRec.Range("A2").Value = "Planilha de pontuação de recursos para o jogo Jogo da WHIP classificada por pontuação média"
$wsRec.Range("A2:H2").HorizontalAlignment = -4131
$wsRec.Range("A2:H2").Borders.Item(9).LineStyle = 1
$wsRec.Range("A2:H2").Borders.Item(9).Weight = 2

# header row (row 3)
$wsRec.Range("A3").Value = "Categoria"
$wsRec.Range("B3").Value = "Recurso"
$wsRec.Range("C3").Value = "Produção"
$wsRec.Range("D3").Value = "Arte"
$wsRec.Range("E3").Value = "Design"
$wsRec.Range("F3").Value = "Engenharia"
$wsRec.Range("G3").Value = "Analista de Qualidade"
$wsRec.Range("H3").Value = "Média"
$wsRec.Range("A3:H3").Font.Bold = $true
$wsRec.Range("A3:H3").Font.ColorIndex = 1
$wsRec.Range("A3:H3").Interior.Pattern = 1
$wsRec.Range("A3:H3").HorizontalAlignment = -4108
$wsRec.Range("A3:H3").VerticalAlignment = -4160
$wsRec.Range("A3:H3").Borders.LineStyle = 1
$wsRec.Range("G3").WrapText = $true
$wsRec.Rows.Item(3).RowHeight = 30

# data rows 4-8 (already in final descending-by-average order)
$data = @(
    @("Processo", "estabelecer um sistema para a circulação de documentos de design e de atualizações de documentos entre a equipe", 3, 3, 3, 3, 3),
    @("Jogabilidade", "interface de usuário fácil de entender", 3, 3, 3, 2, 3),
    @("Produção", "melhorias das questões relacionadas a movimentação do cursor no tabuleiro auxiliar", 2, 1, 2, 2, 2),
    @("Jogabilidade", "possibilidade de o jogador contabilizar as partidas vencidas", 2, 2, 2, 1, 2),
    @("Jogabilidade", "possibilidade dos jogadores personalizarem a aparência dos personagens", 1, 2, 1, 1, 1)
)

$r = 4
foreach ($row in $data) {
    $wsRec.Range("A$r").Value = $row[0]
    $wsRec.Range("B$r").Value = $row[1]
    $wsRec.Range("C$r").Value = $row[2]
    $wsRec.Range("D$r").Value = $row[3]
    $wsRec.Range("E$r").Value = $row[4]
    $wsRec.Range("F$r").Value = $row[5]
    $wsRec.Range("G$r").Value = $row[6]
    $wsRec.Range("H$r").Formula = "=AVERAGE(C$r`:G$r)"
    $r++
}

$wsRec.Range("A4:H8").Borders.LineStyle = 1
$wsRec.Range("A4:B8").HorizontalAlignment = -4131
$wsRec.Range("A4:B8").VerticalAlignment = -4160
$wsRec.Range("B4:B8").WrapText = $true
$wsRec.Rows.Item(4).RowHeight = 30
$wsRec.Rows.Item(6).RowHeight = 30

# legend block (rows 10-12), with a holding column of blank wrap-text cells
for ($i = 9; $i -le 36; $i++) {
    $wsRec.Range("B$i").WrapText = $true
}

$wsRec.Range("A10").Value = 3
$wsRec.Range("B10").Value = "necessário"
$wsRec.Range("A11").Value = 2
$wsRec.Range("B11").Value = "desejado"
$wsRec.Range("A12").Value = 1
$wsRec.Range("B12").Value = "interessante"

$wsRec.Range("G10").Select()

# ---------------------------------------------------------------------------
# 3. Selection / active-tab bookkeeping so the saved file opens the same way
# ---------------------------------------------------------------------------
$wsConceito = $wb.Worksheets.Item("descricao-fase-de-conceituacao")
$wsConceito.Activate()
$wsConceito.Range("E9").Select()

$wsSwot.Activate()
$wsSwot.Range("B10").Select()
